$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Capture the hyperlink target addresses that live on sheet1 (rows 3-6, column A)
# before we touch anything, keyed by source row.
$links = @()
foreach ($hl in $ws1.Hyperlinks) {
    $links += , @($hl.Range.Row, $hl.Address)
}

# Move rows 3:6 (A3:D6) from sheet1 down to sheet2, landing at A1:D4 (values + formats).
$ws1.Range("A3:D6").Copy($ws2.Range("A1"))

# Recreate the hyperlinks on sheet2 (row offset -2: old row 3 -> new row 1, etc.),
# then restore the plain (non-hyperlink) cell style the source cells had - Excel's
# Hyperlinks.Add forces the built-in "Hyperlink" style onto the target cell.
foreach ($pair in $links) {
    $oldRow = $pair[0]
    $addr = $pair[1]
    $newRow = $oldRow - 2
    $target = $ws2.Cells.Item($newRow, 1)
    $ws2.Hyperlinks.Add($target, $addr) | Out-Null
    $ws2.Cells.Item($newRow, 2).Copy()
    $target.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}
$excel.CutCopyMode = $false

# Remove the hyperlinks that used to live on sheet1 (now duplicated on sheet2) and
# delete the rows that were moved off of sheet1.
$ws1.Hyperlinks.Delete()
$ws1.Rows("3:6").Delete()

# Update sheet2's selection to span the whole used range (A1:XFD4), as in the target.
$ws2.Range("A1:XFD4").Select() | Out-Null

# Restore sheet1's original selection/active cell (A2) and make it the active sheet again.
$ws1.Activate() | Out-Null
$ws1.Range("A2").Select() | Out-Null
